$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 228.07692
$ws.Range("I11").Value = 228.07692
$ws.Range("K11").Value = 228.07692
$ws.Range("M11").Value = -88.07692
$ws.Range("H43").Value = 1887.25
$ws.Range("I43").Value = 1849.5
$ws.Range("J43").Value = 1925
$ws.Range("K43").Value = 1849.5
$ws.Range("L43").Value = 1925
$ws.Range("M43").Value = -1780.5
$ws.Range("N43").Value = -2063
$ws.Range("H103").Value = 2117.2222
$ws.Range("I103").Value = 2579.2856
$ws.Range("J103").Value = 500
$ws.Range("K103").Value = 7737.8568
$ws.Range("L103").Value = 1500
$ws.Range("M103").Value = -7151.8568
$ws.Range("N103").Value = -2672
$ws.Range("H136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("N136").Value = 0
$ws.Range("H138").Value = 3138.6904
$ws.Range("J138").Value = 3753.7144
$ws.Range("L138").Value = 11261.1432
$ws.Range("N138").Value = -21541.1432
$ws.Range("H141").Value = 1916
$ws.Range("I141").Value = 1816.762
$ws.Range("K141").Value = 5450.286
$ws.Range("M141").Value = -270.2860000000001

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 3873.5
$ws.Range("I63").Value = 2498.3333
$ws.Range("K63").Value = 2498.3333
$ws.Range("M63").Value = -1812.3333
$ws.Range("H66").Value = 3873.5
$ws.Range("I66").Value = 2498.3333
$ws.Range("K66").Value = 12491.6665
$ws.Range("M66").Value = -9059.666499999999
$ws.Range("H74").Value = 243573.88
$ws.Range("J74").Value = 2999.8333
$ws.Range("L74").Value = 2999.8333
$ws.Range("N74").Value = -4747.8333
$ws.Range("H77").Value = 243573.88
$ws.Range("J77").Value = 2999.8333
$ws.Range("L77").Value = 14999.1665
$ws.Range("N77").Value = -23735.1665
$ws.Range("H97").Value = 1388.3478
$ws.Range("I97").Value = 1316.8948
$ws.Range("K97").Value = 1316.8948
$ws.Range("M97").Value = -820.8948
$ws.Range("H132").Value = 3072.9333
$ws.Range("I132").Value = 1998.5312
$ws.Range("K132").Value = 5995.5936
$ws.Range("M132").Value = -3465.5936

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 76928670
$ws.Range("I94").Value = 105269810
$ws.Range("K94").Value = 105269810
$ws.Range("M94").Value = -105269359
$ws.Range("H105").Value = 12382546
$ws.Range("I105").Value = 1001343.6
$ws.Range("K105").Value = 1001343.6
$ws.Range("M105").Value = -999596.6

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1478.4736
$ws.Range("I16").Value = 1699.1538
$ws.Range("K16").Value = 1699.1538
$ws.Range("M16").Value = -1412.1538
$ws.Range("H31").Value = 1897565
$ws.Range("I31").Value = 2934.0625
$ws.Range("K31").Value = 2934.0625
$ws.Range("M31").Value = -2639.0625
$ws.Range("H34").Value = 1897565
$ws.Range("I34").Value = 2934.0625
$ws.Range("K34").Value = 2934.0625
$ws.Range("M34").Value = -2732.0625
$ws.Range("H41").Value = 18840.5
$ws.Range("I41").Value = 2939.75
$ws.Range("J41").Value = 34741.25
$ws.Range("K41").Value = 2939.75
$ws.Range("L41").Value = 34741.25
$ws.Range("M41").Value = -2511.75
$ws.Range("N41").Value = -35597.25
$ws.Range("H58").Value = 2662.3684
$ws.Range("I58").Value = 2389.577
$ws.Range("J58").Value = 3253.4167
$ws.Range("K58").Value = 2389.577
$ws.Range("L58").Value = 3253.4167
$ws.Range("M58").Value = -2186.577
$ws.Range("N58").Value = -3659.4167
$ws.Range("H62").Value = 14289355
$ws.Range("I62").Value = 16670248
$ws.Range("K62").Value = 16670248
$ws.Range("M62").Value = -16669624
$ws.Range("H65").Value = 14289355
$ws.Range("I65").Value = 16670248
$ws.Range("K65").Value = 83351240
$ws.Range("M65").Value = -83348120
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("N99").Value = 0
$ws.Range("H105").Value = 1014
$ws.Range("I105").Value = 1014
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 1014
$ws.Range("L105").Value = 0
$ws.Range("N105").Value = 733
$ws.Range("H113").Value = 1478.4736
$ws.Range("I113").Value = 1699.1538
$ws.Range("K113").Value = 1699.1538
$ws.Range("M113").Value = 470.8462
$ws.Range("H122").Value = 3233.7693
$ws.Range("I122").Value = 3008.889
$ws.Range("K122").Value = 9026.667000000001
$ws.Range("M122").Value = -6576.667000000001
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("N126").Value = 0
$ws.Range("H132").Value = 22227730
$ws.Range("I132").Value = 2179.2222
$ws.Range("K132").Value = 6537.6666
$ws.Range("M132").Value = -4007.6666
$ws.Range("H136").Value = 2662.3684
$ws.Range("I136").Value = 2389.577
$ws.Range("J136").Value = 3253.4167
$ws.Range("K136").Value = 7168.731000000001
$ws.Range("L136").Value = 9760.250100000001
$ws.Range("M136").Value = -4618.731000000001
$ws.Range("N136").Value = -14860.2501

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("M70").ClearContents()
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("M73").ClearContents()
$ws.Range("H80").Value = 200
$ws.Range("J80").Value = 200
$ws.Range("L80").Value = 600
$ws.Range("N80").Value = -2472
$ws.Range("H83").Value = 200
$ws.Range("J83").Value = 200
$ws.Range("L83").Value = 1800
$ws.Range("N83").Value = -11160

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 14579.7
$ws.Range("I126").Value = 4059.6
$ws.Range("K126").Value = 12178.8
$ws.Range("M126").Value = -9708.799999999999

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1934.75
$ws.Range("I22").Value = 444
$ws.Range("J22").Value = 2431.6667
$ws.Range("K22").Value = 444
$ws.Range("L22").Value = 2431.6667
$ws.Range("M22").Value = -149
$ws.Range("N22").Value = -3021.6667
$ws.Range("H27").Value = 1934.75
$ws.Range("I27").Value = 444
$ws.Range("J27").Value = 2431.6667
$ws.Range("K27").Value = 444
$ws.Range("L27").Value = 2431.6667
$ws.Range("M27").Value = -337
$ws.Range("N27").Value = -2645.6667
$ws.Range("H45").Value = 30041
$ws.Range("I45").Value = 30041
$ws.Range("K45").Value = 30041
$ws.Range("M45").Value = -29634
$ws.Range("H46").Value = 1513.6
$ws.Range("I46").Value = 1517
$ws.Range("K46").Value = 1517
$ws.Range("M46").Value = -1329
$ws.Range("H55").Value = 663
$ws.Range("I55").Value = 372.75
$ws.Range("J55").Value = 1050
$ws.Range("K55").Value = 372.75
$ws.Range("L55").Value = 1050
$ws.Range("M55").Value = -199.75
$ws.Range("N55").Value = -1396
$ws.Range("H122").Value = 5369.4443
$ws.Range("I122").Value = 3243.4666
$ws.Range("K122").Value = 9730.399800000001
$ws.Range("M122").Value = -7280.399800000001
$ws.Range("H132").Value = 3166.8
$ws.Range("I132").Value = 3166.8
$ws.Range("K132").Value = 9500.400000000001
$ws.Range("M132").Value = -6970.400000000001

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H7").Value = 2004
$ws.Range("I7").Value = 2004
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 2004
$ws.Range("M7").Value = -1891
$ws.Range("N7").Value = 0
$ws.Range("H9").Value = 153372.75
$ws.Range("J9").Value = 153372.75
$ws.Range("L9").Value = 153372.75
$ws.Range("N9").Value = -153652.75
$ws.Range("H10").Value = 9333.333000000001
$ws.Range("J10").Value = 9333.333000000001
$ws.Range("L10").Value = 9333.333000000001
$ws.Range("N10").Value = -9671.333000000001
$ws.Range("H14").Value = 10004
$ws.Range("J14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("N14").ClearContents()
$ws.Range("H17").Value = 7500
$ws.Range("J17").Value = 7333.3335
$ws.Range("L17").Value = 7333.3335
$ws.Range("N17").Value = -7677.3335
$ws.Range("H62").Value = 6790.3335
$ws.Range("J62").Value = 8078.6
$ws.Range("L62").Value = 8078.6
$ws.Range("N62").Value = -9326.6
$ws.Range("H65").Value = 6790.3335
$ws.Range("J65").Value = 8078.6
$ws.Range("L65").Value = 40393
$ws.Range("N65").Value = -46633
$ws.Range("H95").Value = 43333.332
$ws.Range("I95").Value = 43000
$ws.Range("K95").Value = 43000
$ws.Range("M95").Value = -40254
$ws.Range("H132").Value = 2847.724
$ws.Range("I132").Value = 2903.48
$ws.Range("K132").Value = 8710.440000000001
$ws.Range("M132").Value = -6180.440000000001
$ws.Range("H135").Value = 64954.5
$ws.Range("J135").Value = 64954.5
$ws.Range("L135").Value = 64954.5
$ws.Range("N135").Value = -75094.5
